# Fruta / hortaliza, semanal
# Update the weekly Níspero price records (rows 3-9) to reflect the new
# weekly data pull. Only the cells whose values actually changed are
# touched; column headers (row 1) and the rest of each record stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44902
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/caja 10 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1500

# Row 4
$ws.Range("D4").Value = 44902
$ws.Range("M4").Value = 70
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 13000
$ws.Range("Q4").Value = "$/caja 10 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1300

# Row 5
$ws.Range("D5").Value = 45251
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 150
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("Q5").Value = "$/bandeja 10 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 2000

# Row 6
$ws.Range("D6").Value = 44505
$ws.Range("K6").Value = "Californiana(o)"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("Q6").Value = "$/bandeja 10 kilos"
$ws.Range("S6").Value = 1500

# Row 7
$ws.Range("D7").Value = 44505
$ws.Range("M7").Value = 50
$ws.Range("Q7").Value = "$/bandeja 10 kilos"

# Row 8
$ws.Range("D8").Value = 45264
$ws.Range("K8").Value = "Golden Nugget"
$ws.Range("M8").Value = 150
$ws.Range("Q8").Value = "$/caja 10 kilos"

# Row 9
$ws.Range("D9").Value = 45250
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 2000
